$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1: "Gender(M/F)" with "Gender(M/F" bold and ")" regular.
$ws.Range("G1").Value = "Gender(M/F)"

$bold = $ws.Range("G1").Characters(1, 10)
$bold.Font.Bold = $true
$bold.Font.Name = "Calibri"
$bold.Font.Size = 11

$rest = $ws.Range("G1").Characters(11, 1)
$rest.Font.Bold = $false
$rest.Font.Name = "Calibri"
$rest.Font.Size = 11

# Widen the new column like the rest of the header columns.
$ws.Columns.Item(7).ColumnWidth = 12.42578125

# Selection collapses back to a single cell after the edit.
$ws.Range("A2").Select()
